$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Title text: "Assignment 1" -> "Assignment 2", and leave a
#    "_GoBack" bookmark at the edit point (mirrors what real Word
#    does after you type over a selection).
# ------------------------------------------------------------------
$titleRange = $d.Content
$titleRange.Find.Execute("Assignment 1", $false, $false, $false, $false, $false, `
                          $true, 1, $false, "Assignment 2", 2) | Out-Null

$goBackRange = $d.Range($titleRange.End - 1, $titleRange.End)
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null

# ------------------------------------------------------------------
# 2. Re-create the seven legacy TOC bookmarks that wrap
#    "Analysis and Design Document" in their original document order
#    (deleting and re-adding makes the in-package bookmark ids come
#    out renumbered the way the reference edit shows).
# ------------------------------------------------------------------
$headingNames = @(
    "_Toc254785383",
    "_Toc254771757",
    "_Toc254770266",
    "_Toc254770226",
    "_Toc222883075",
    "_Toc222821167",
    "_Toc222820221"
)
foreach ($n in $headingNames) {
    if ($d.Bookmarks.Exists($n)) {
        $d.Bookmarks($n).Delete()
    }
}

$headingRange = $d.Content
$headingRange.Find.Execute("Analysis and Design Document") | Out-Null

$addOrder = @(
    "_Toc222820221",
    "_Toc222821167",
    "_Toc222883075",
    "_Toc254770226",
    "_Toc254770266",
    "_Toc254771757",
    "_Toc254785383"
)
foreach ($n in $addOrder) {
    $d.Bookmarks.Add($n, $headingRange) | Out-Null
}
